$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update "Version" value (row 3, column B): 0.1.1 -> 0.2.0 ---
$ws.Range("B3").Value = "0.2.0"

# --- Update "Date" value (row 8, column B) ---
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# --- Insert a new "Jurisdiction" / "iso:code:3166:FR" row after the
#     "Contact" row (row 10), pushing the Description/Purpose/Copyright/
#     Immutable rows down by one (old rows 11-14 -> new rows 12-15). ---

# Capture the existing (pre-shift) text of rows 11-14 before overwriting.
$a11 = $ws.Range("A11").Text
$b11 = $ws.Range("B11").Text
$a12 = $ws.Range("A12").Text
$b12 = $ws.Range("B12").Text
$a13 = $ws.Range("A13").Text
$b13 = $ws.Range("B13").Text
$a14 = $ws.Range("A14").Text
$b14 = $ws.Range("B14").Text

# Extend formatting (border/alignment) from row 14 down onto the new row 15.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Shift content down one row at a time, from the bottom up.
$ws.Range("A15").Value = $a14
$ws.Range("B15").Value = $b14

$ws.Range("A14").Value = $a13
$ws.Range("B14").Value = $b13

$ws.Range("A13").Value = $a12
$ws.Range("B13").Value = $b12

$ws.Range("A12").Value = $a11
$ws.Range("B12").Value = $b11

# Write the new Jurisdiction row into the now-vacated row 11.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
